$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()
$ws.Range("A7").Value = "This is in staging area"
$ws.Range("A7").Select()
